# Estadisticos Matutinos 15 Oct
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 1P" (sheet1) - update stats for the 4 groups
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Range("D2").Value = 17
$ws1.Range("F2").Value = 18
$ws1.Range("G2").Value = 51.43
$ws1.Range("H2").Value = 6.8

$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 5
$ws1.Range("F3").Value = 34
$ws1.Range("G3").Value = 87.18000000000001
$ws1.Range("H3").Value = 7.7

$ws1.Range("D4").Value = 5
$ws1.Range("F4").Value = 30
$ws1.Range("G4").Value = 85.70999999999999
$ws1.Range("H4").Value = 8.1

$ws1.Range("D5").Value = 17
$ws1.Range("F5").Value = 16
$ws1.Range("G5").Value = 48.48
$ws1.Range("H5").Value = 7.2

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P" (sheet2) - only the "Reprobados" column changes
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Range("E2").Value = 18
$ws2.Range("E3").Value = 39
$ws2.Range("E4").Value = 30
$ws2.Range("E5").Value = 16

# ---------------------------------------------------------------------------
# Sheet "Estadisticos Final" (sheet3) - same update as sheet1
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Range("D2").Value = 17
$ws3.Range("F2").Value = 18
$ws3.Range("G2").Value = 51.43
$ws3.Range("H2").Value = 6.8

$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 5
$ws3.Range("F3").Value = 34
$ws3.Range("G3").Value = 87.18000000000001
$ws3.Range("H3").Value = 7.7

$ws3.Range("D4").Value = 5
$ws3.Range("F4").Value = 30
$ws3.Range("G4").Value = 85.70999999999999
$ws3.Range("H4").Value = 8.1

$ws3.Range("D5").Value = 17
$ws3.Range("F5").Value = 16
$ws3.Range("G5").Value = 48.48
$ws3.Range("H5").Value = 7.2

# ---------------------------------------------------------------------------
# Sheet "Rescatables" (sheet4) - add the 8 new students (rows 2-9)
# Values are entered column by column so new shared strings are created in
# the same order as in the target workbook (Paterno, then Materno, then
# Nombres columns).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$matriculas = @(20330051920072, 20330051920073, 20330051920082, 20330051920070, 20330051920158, 20330051920306, 20330051920090, 20330051920226)
$paterno    = @("CARRERA", "CANUTO", "GONZALEZ", "TEXCAHUA", "BAUTISTA", "PARRA", "MORALES", "GARCIA")
$materno    = @("CASTAÑEDA", "MEDINA", "OFICIAL", "MARTINEZ", "DIAZ", "FLORES", "BAROJAS", "CRUZ")
$nombres    = @("PAUL ARAVIER", "ISRAEL", "SAID ANDRES", "OSVALDO", "DINA BERENICE", "SUEMI", "DIEGO IVAN", "JESSICA")
$carreras   = @("BIOLOGÍA", "BIOLOGÍA", "BIOLOGÍA", "BIOLOGÍA", "BIOLOGÍA", "BIOLOGÍA", "BIOLOGÍA", "BIOLOGÍA")
$grupos     = @("3AEV", "3AEV", "3AEV", "3AEV", "3ARHM", "3ARHM", "3AEV", "3ALCM")
$reprobadas = @(6, 6, 6, 6, 6, 6, 6, 6)

for ($i = 0; $i -lt $matriculas.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 1).Value = $matriculas[$i]
}
for ($i = 0; $i -lt $paterno.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt $materno.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt $nombres.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt $carreras.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 5).Value = $carreras[$i]
}
for ($i = 0; $i -lt $grupos.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 6).Value = $grupos[$i]
}
for ($i = 0; $i -lt $reprobadas.Length; $i++) {
    $row = $i + 2
    $ws4.Cells.Item($row, 7).Value = $reprobadas[$i]
}
